$d = $word.ActiveDocument

# --- Locate the paragraph that currently ends the document's body content ---
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*otherwise difficult to reproduce.*") {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    throw "Could not find the anchor paragraph ending in 'otherwise difficult to reproduce.'"
}

$insertAt = $target.Range.End
$insertionPoint = $d.Range($insertAt, $insertAt)

# --- Build the two new paragraphs as a raw WordprocessingML fragment so the
#     exact run/proofErr structure from the source edit is reproduced. The
#     "pingall" run is left as plain text here; its Emphasis character style
#     is applied afterwards via Range.Style (InsertXML silently drops
#     w:rStyle references). ---
$newParagraphs = @"
<w:p>
<w:pPr>
<w:pStyle w:val="Heading2"/>
</w:pPr>
<w:r>
<w:t>Customizing the Controller</w:t>
</w:r>
</w:p>
<w:p>
<w:r>
<w:tab/>
<w:t xml:space="preserve">The controller is responsible for making all flow decisions as the data plane is simple packet forwarding devices. </w:t>
</w:r>
<w:proofErr w:type="spellStart"/>
<w:r>
<w:t>Dordal</w:t>
</w:r>
<w:proofErr w:type="spellEnd"/>
<w:r>
<w:t xml:space="preserve"> demonstrates this by constructing a rectangular looped topology, then launching it without a controller. After issuing </w:t>
</w:r>
<w:proofErr w:type="spellStart"/>
<w:r>
<w:t>mininet’s</w:t>
</w:r>
<w:proofErr w:type="spellEnd"/>
<w:r>
<w:t xml:space="preserve"> </w:t>
</w:r>
<w:proofErr w:type="spellStart"/>
<w:r>
<w:t>pingall</w:t>
</w:r>
<w:proofErr w:type="spellEnd"/>
<w:r>
<w:t xml:space="preserve"> test command, an infinite loop occurs.</w:t>
</w:r>
</w:p>
"@

$xmlPackage = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $newParagraphs + '</w:body></w:document>' +
  '</pkg:xmlData>' +
  '</pkg:part>' +
  '</pkg:package>'

$insertionPoint.InsertXML($xmlPackage)

# --- Apply the "Emphasis" character style to the "pingall" run in the newly
#     inserted paragraph. ---
$searchRange = $d.Range($insertAt, $d.Content.End)
$found = $searchRange.Find.Execute("pingall", $true, $false, $false, $false, $false,
                                    $true, 1, $false, "", 0)
if ($found) {
    $searchRange.Style = "Emphasis"
}

# --- The "_GoBack" bookmark must end up right after the final run of the
#     newly-added content (i.e. at the very end of the document), matching
#     where Word leaves it after the last edit. Re-adding a bookmark with the
#     same name moves/replaces the existing one. ---
$docEnd = $d.Content.End
$bookmarkRange = $d.Range($docEnd, $docEnd)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
